$wb = $excel.ActiveWorkbook

# This handoff report regenerates the status of the 614fb894 file: it has
# moved from "Handed back: in sync with en-US" to "In Translation", its
# handoff timestamps were refreshed, and an "out of date handback" error
# detail was attached on both locale sheets.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("G3").Value = "2017-02-17 09:38:05"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("H3").Value = "2017-02-17 09:37:49"
$zhcn.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4628176429fc25ef3f5c5e2690b4744eb70185c/e2e/614fb894-1524-4b63-a580-394b1a45dc50.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb49a54b6e461f7f27fa510be098785493daea29/e2e/614fb894-1524-4b63-a580-394b1a45dc50.md."

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "In Translation"
$dede.Range("H3").Value = "2017-02-17 09:38:05"
$dede.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4628176429fc25ef3f5c5e2690b4744eb70185c/e2e/614fb894-1524-4b63-a580-394b1a45dc50.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb49a54b6e461f7f27fa510be098785493daea29/e2e/614fb894-1524-4b63-a580-394b1a45dc50.md."
